$wb = $excel.ActiveWorkbook

# ---- Metrics sheet: update source values ----
$metrics = $wb.Worksheets.Item("Metrics")

$metrics.Range("B2").Value = 335196.5
$metrics.Range("B3").Value = 274840.93000000005
$metrics.Range("B4").Value = 106817.01999999999
$metrics.Range("B5").Value = 13385
$metrics.Range("B6").Value = 4702327.97
$metrics.Range("B7").Value = 3964659.5999999992
$metrics.Range("B8").Value = 1377419.16
$metrics.Range("B9").Value = 182386
$metrics.Range("B10").Value = 33167651.770999823
$metrics.Range("B11").Value = 31239881.120000001
$metrics.Range("B12").Value = 11659128.050000003
$metrics.Range("B13").Value = 1280013

# Move the saved selection on the Metrics sheet
$metrics.Range("E24").Select() | Out-Null

# ---- today sheet: move the saved selection (values recalc automatically) ----
$today = $wb.Worksheets.Item("today")
$today.Activate() | Out-Null
$today.Range("F7").Select() | Out-Null
